# edit.ps1 -- "updated spec and personas"
#
# Turns each one-paragraph persona blurb into a two-paragraph pair:
#   "<Name>, <age> years old"
#   "As a <role> ... So that <benefit>."   (an agile user-story)
# and tidies up the blank paragraphs that separated/followed them.
#
# Personas are located by searching for their distinctive original text
# (rather than hard-coded paragraph indices) so the script keeps working
# even if the surrounding structure shifts a little, and the three
# replacements are done back-to-front (Jack, then Jill, then John) so that
# each `Find` below still locates the right, not-yet-edited paragraph.

$d = $word.ActiveDocument

function Find-ParagraphStartingWith($doc, $needle) {
    $r = $doc.Content
    $ok = $r.Find.Execute($needle, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        throw "Could not find a paragraph starting with: $needle"
    }
    return $r.Paragraphs(1)
}

# --- Jack: "Jack, Teacher at high school ..." -> name paragraph + user story (keeps the _GoBack bookmark) ---
$jackXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>Jack,</w:t></w:r><w:r><w:t xml:space="preserve"> 46 years old</w:t></w:r></w:p><w:p><w:r><w:t>As a</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>t</w:t></w:r><w:r><w:t xml:space="preserve">eacher at high school </w:t></w:r><w:r><w:t>I want an interesting graphical display of non-communicable diseases so that I can teach my students about them and help them be more informed. I am particularly interested in clear graphic displays where students could take in data quickly.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$pJack = Find-ParagraphStartingWith $d "Jack, Teacher"
$pJack.Range.InsertXML($jackXml)

# --- Jill: "Jill, regular person ..." -> name paragraph + user story ---
$jillXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>Jill,</w:t></w:r><w:r><w:t xml:space="preserve"> 38 years old</w:t></w:r></w:p><w:p><w:r><w:t>As a civilian</w:t></w:r><w:r><w:t xml:space="preserve"> with </w:t></w:r><w:r><w:t xml:space="preserve">an </w:t></w:r><w:r><w:t xml:space="preserve">interest in the symptoms of non-communicable diseases. </w:t></w:r><w:r><w:t xml:space="preserve">I want to look up diseases and symptoms about them to be more informed on them. </w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>So that I can assist with my family member effected by the disease.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$pJill = Find-ParagraphStartingWith $d "Jill, regular person"
$pJill.Range.InsertXML($jillXml)

# --- John: "John, Student at university ..." -> name paragraph + user story ---
$johnXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>John, 20 years old</w:t></w:r></w:p><w:p><w:r><w:t>As a</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>s</w:t></w:r><w:r><w:t xml:space="preserve">tudent at university who has a project based on non-communicable </w:t></w:r><w:r><w:t>diseases, I want to find</w:t></w:r><w:r><w:t xml:space="preserve"> information sources about symptoms and trends in data</w:t></w:r><w:r><w:t>. So that I can make a great project and achieve a high grade.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$pJohn = Find-ParagraphStartingWith $d "John, Student"
$pJohn.Range.InsertXML($johnXml)

# --- remove the blank paragraph right after the "User Personas" title ---
$pTitle = Find-ParagraphStartingWith $d "User Personas"
$pBlankAfterTitle = $pTitle.Next()
$pBlankAfterTitle.Range.Delete()

# --- remove the blank paragraph that used to trail Jack's blurb, at the very end of the document ---
# The very last paragraph mark of a document can't be selected/deleted on
# its own, so instead delete the mark that ends the paragraph just before
# it; that merges the (now empty) trailing paragraph away.
$pLast = $d.Paragraphs($d.Paragraphs.Count)
$pBeforeLast = $pLast.Previous()
$endMark = $d.Range($pBeforeLast.Range.End - 1, $pBeforeLast.Range.End)
$endMark.Delete()
